$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from SCD0252 to SCD0016
$ws.Name = "SCD0016"

# Update cell B2's value from "DGS-267" to "SCD0016-026"
$ws.Range("B2").Value = "SCD0016-026"

# Apply the same font as C2/D2 (Arial 10) to B2 so it matches the updated style
$ws.Range("B2").Font.Name = $ws.Range("C2").Font.Name
$ws.Range("B2").Font.Size = $ws.Range("C2").Font.Size

# Widen column B to fit the new value
$ws.Columns("B").ColumnWidth = 11.67

# Move the active selection to B2
$ws.Range("B2").Select() | Out-Null
